# Generate Report for Handoff
# Updates the localization-status report to reflect that the zh-cn content
# is now "Ready for handoff" (was "In Translation"), refreshes the relevant
# timestamps, and widens the status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-13 05:11:57"

# --- zh-cn detail sheet -----------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-13 05:11:50"

# --- de-de detail sheet -----------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime (shares the generate-date string)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-13 05:11:57"

# --- Widen the Status columns to fit "Ready for handoff" --------------------
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
